$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Boolean" sheet: split the combined BVTQaZ and VTQaZ CSV rows into the new
# per-vehicle-type CSV rows, and leave a handful of blank rows at the bottom
# (mirrors what the author did in Excel: insert rows, type new paths).
# ---------------------------------------------------------------------------
$boolSheet = $wb.Worksheets.Item("Boolean")

# Row 17 currently holds "trans/BVTQaZ/BVTQaZ.csv" - insert 5 rows below it
# so there is room for 6 total rows, then fill them in.
$boolSheet.Range("A18:A22").Insert()
$boolSheet.Range("A17").Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$boolSheet.Range("A18").Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$boolSheet.Range("A19").Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$boolSheet.Range("A20").Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$boolSheet.Range("A21").Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$boolSheet.Range("A22").Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# Rows 23-25 are now "trans/BVTStL/BVTStL.csv", "trans/PVTStL/PVTStL.csv" and
# "trans/SRPbVT/SRPbVT.csv" - untouched. Row 26 holds "trans/VTQaZ/VTQaZ.csv";
# insert 5 more rows below it for the same split.
$boolSheet.Range("A27:A31").Insert()
$boolSheet.Range("A26").Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$boolSheet.Range("A27").Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$boolSheet.Range("A28").Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$boolSheet.Range("A29").Value = "trans/VTQaZ/VTQaZ-rail.csv"
$boolSheet.Range("A30").Value = "trans/VTQaZ/VTQaZ-ships.csv"
$boolSheet.Range("A31").Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# Row 32 is "trans/VTStFES/VTStFES.csv" (unchanged). Add six blank rows after
# it, formatted like the rest of the list (style carried by column A).
$boolSheet.Range("A33:A38").Value = ""

$boolSheet.Range("A17:A38").Font.Name = "Calibri"

# ---------------------------------------------------------------------------
# View/selection bookkeeping to mirror the saved workbook state: the "About"
# sheet becomes the active tab, and the last-used cell on "Integer" /
# "Boolean" move to where the author was last working.
# ---------------------------------------------------------------------------
$integerSheet = $wb.Worksheets.Item("Integer")
$integerSheet.Activate()
$integerSheet.Range("A13").Select()

$boolSheet.Activate()
$boolSheet.Range("A32").Select()

$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Activate()
$aboutSheet.Range("A1").Select()
